$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts the existing A:E
# columns (segment names, PercActivations, PercSegmentAreas,
# RelativeCAMImportance, PercActivationsRescaled) one column to the
# right, becoming B:F, carrying their formatting along with them.
# The header row had no A1 cell, so after the shift B1 is empty too.
$ws.Columns.Item(1).Insert()

# New header, in the now-empty B1, matching the formatting already
# used by the rest of row 1 (bold, bordered, centered).
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B1").Value = "segments"

# Fill the new column A with the 0-based segment index for the 19
# data rows (rows 2-20), matching the style previously used for the
# segment-name column (bold, bordered, centered) -- copy it *before*
# stripping that same styling from the (now plain-text) segment
# names in column B.
$ws.Range("B2").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt 19; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
}

# The segment-name column (now B) inherited the bold/bordered header
# styling it used to carry in column A; the new layout wants that
# column to be plain/unstyled text.
$ws.Range("B2:B20").ClearFormats()
